# Natmi following Dr Hou advice
#
# The original sheet held the Cthrc1 -> Ror2 ligand/receptor table for a
# single sending/target cluster pair (FAPs -> FAPs / sCs). Per Dr Hou's
# advice, the cell-type universe was expanded to include "ECs" (sending)
# and "Neutro" (target), and every sending/target cluster combination is
# now reported (3 sending clusters x 3 target clusters = 9 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each inner array is one data row, in column order A..T:
#   A: Sending cluster, B: Ligand symbol, C: Receptor symbol, D: Target cluster,
#   E..T: the numeric NATMI statistics for that sending/target pair.
$rows = @(
        @("ECs", "Cthrc1", "Ror2", "FAPs", 2, 0.6666666666666666, 0.9744206666666667, 2.923262, 0.0323957296480766, 0.03239572964807661, 3, 1, 4.669265333333333, 14.007796, 0.8833113458668934, 0.8833113458668933, 4.549828638950222, 40.948457750552, 0.02861551555578257, 0.02861551555578257),
        @("ECs", "Cthrc1", "Ror2", "Neutro", 2, 0.6666666666666666, 0.9744206666666667, 2.923262, 0.0323957296480766, 0.03239572964807661, 1, 0.3333333333333333, 0.017474, 0.052422, 0.003305655463074583, 0.003305655463074583, 0.01702702672933333, 0.153243240564, 0.0001070891206914517, 0.0001070891206914517),
        @("ECs", "Cthrc1", "Ror2", "sCs", 2, 0.6666666666666666, 0.9744206666666667, 2.923262, 0.0323957296480766, 0.03239572964807661, 3, 1, 0.5993529999999999, 1.798059, 0.1133829986700321, 0.1133829986700321, 0.5840219498286666, 5.256197548458, 0.003673124971602588, 0.003673124971602589),
        @("FAPs", "Cthrc1", "Ror2", "FAPs", 3, 1, 28.021538, 84.06461399999999, 0.9316080830640275, 0.9316080830640276, 3, 1, 4.669265333333333, 14.007796, 0.8833113458668934, 0.8833113458668933, 130.8399959700826, 1177.559963730744, 0.8228999896717628, 0.8228999896717628),
        @("FAPs", "Cthrc1", "Ror2", "Neutro", 3, 1, 28.021538, 84.06461399999999, 0.9316080830640275, 0.9316080830640276, 1, 0.3333333333333333, 0.017474, 0.052422, 0.003305655463074583, 0.003305655463074583, 0.4896483550119999, 4.406835195108, 0.003079575349225042, 0.003079575349225043),
        @("FAPs", "Cthrc1", "Ror2", "sCs", 3, 1, 28.021538, 84.06461399999999, 0.9316080830640275, 0.9316080830640276, 3, 1, 0.5993529999999999, 1.798059, 0.1133829986700321, 0.1133829986700321, 16.79479286491399, 151.153135784226, 0.1056285180430397, 0.1056285180430398),
        @("sCs", "Cthrc1", "Ror2", "FAPs", 3, 1, 1.082717666666667, 3.248153, 0.03599618728789584, 0.03599618728789584, 3, 1, 4.669265333333333, 14.007796, 0.8833113458668934, 0.8833113458668933, 5.055496066754221, 45.499464600788, 0.03179584063934803, 0.03179584063934803),
        @("sCs", "Cthrc1", "Ror2", "Neutro", 3, 1, 1.082717666666667, 3.248153, 0.03599618728789584, 0.03599618728789584, 1, 0.3333333333333333, 0.017474, 0.052422, 0.003305655463074583, 0.003305655463074583, 0.01891940850733333, 0.170274676566, 0.0001189909931580887, 0.0001189909931580887),
        @("sCs", "Cthrc1", "Ror2", "sCs", 3, 1, 1.082717666666667, 3.248153, 0.03599618728789584, 0.03599618728789584, 3, 1, 0.5993529999999999, 1.798059, 0.1133829986700321, 0.1133829986700321, 0.6489300816696666, 5.840370735026999, 0.004081355655389718, 0.004081355655389718)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowValues = $rows[$i]
    $r = $startRow + $i
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowValues[$col - 1]
    }
}
